# PCB Design Checklist - finish checklist for N.1 board (close #18)
#
# - Update the "Variant" title in A3 (Charge -> IIP_Charge)
# - Flip the Documentation checklist rows 65-72 from NOK ("X" in column D)
#   to OK ("X" in column C) -- row 67 ends up N/A ("X" in column E) instead
# - Add a comment to row 74 (F74)
# - Fix the spelling of "Mechanical validation" in C80
# - Move the selection/scroll position to the bottom of the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")

# --- Row 3: update board variant description ---
$ws.Range("A3").Value = 'Variant: Semi USB - "IIP_Charge" Board Nº1'

# --- Rows 65-66: NOK -> OK ---
$ws.Range("C65").Value = "X"
$ws.Range("D65").ClearContents()

$ws.Range("C66").Value = "X"
$ws.Range("D66").ClearContents()

# --- Row 67: NOK -> N/A ---
$ws.Range("D67").ClearContents()
$ws.Range("E67").Value = "X"

# --- Rows 68-72: NOK -> OK ---
$ws.Range("C68").Value = "X"
$ws.Range("D68").ClearContents()

$ws.Range("C69").Value = "X"
$ws.Range("D69").ClearContents()

$ws.Range("C70").Value = "X"
$ws.Range("D70").ClearContents()

$ws.Range("C71").Value = "X"
$ws.Range("D71").ClearContents()

$ws.Range("C72").Value = "X"
$ws.Range("D72").ClearContents()

# --- Row 74: add comment ---
$ws.Range("F74").Value = "Under development"

# --- Row 80: fix spelling ---
$ws.Range("C80").Value = "Mechanical validation"

# --- Update the active selection / scroll position ---
$ws.Activate()
$ws.Range("C81").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
